# Regenerate orders with updated distance/size labels.
# Simple, consistent substring renames across every text cell in the sheet:
#   D64 -> D69
#   D80 -> D86
#   D51 -> D55
#   S30 -> S31
# (applied in this order so none of the replacement tokens collide with
#  any of the other source tokens)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count

for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value()
        if ($val -ne $null -and $val.GetType().Name -eq "String") {
            if ($val.Contains("D64") -or $val.Contains("D80") -or $val.Contains("D51") -or $val.Contains("S30")) {
                $newVal = $val.Replace("D64", "D69").Replace("D80", "D86").Replace("D51", "D55").Replace("S30", "S31")
                $cell.Value = $newVal
            }
        }
    }
}
